$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the Chauchard, Klasnja, and Harris "Survey" record that used
# digitally-measured chart CIs). This shifts the old rows 14 (Solaz et al.)
# and 15 (Arvate & Mittlaender) up to become rows 13 and 14, and also drops
# the now-unused shared strings describing that record.
$ws.Rows.Item(13).Delete()
